$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns -- and, for the
# two rows that swapped rank this run, the "Coin"/"Link" (B/C) columns
# too -- with the latest scrape. All of these columns are plain text in
# the source sheet (even the numeric-looking "Price" values), so for any
# value that Excel would otherwise auto-convert to a Number we write it
# with a leading apostrophe to force text, then snap the cell's style back
# to Normal so no stray "quote prefix" formatting lingers on the cell.
$ws.Range('D2').Value = '28.499.21'
$ws.Range('E2').Value = '  +5.06%  '
$ws.Range('D3').Value = '1.601.75'
$ws.Range('E3').Value = '  +2.81%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '''215.00'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.53%  '
$ws.Range('E6').Value = '  +1.85%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '''24.02'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +9.20%  '
$ws.Range('E9').Value = '  +1.61%  '
$ws.Range('E10').Value = '  +1.20%  '
$ws.Range('D11').Value = '''0.0890'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.39%  '
$ws.Range('D12').Value = '1.831.15'
$ws.Range('E12').Value = '  +2.82%  '
$ws.Range('D13').Value = '1.609.76'
$ws.Range('E13').Value = '  +3.24%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = '''0.534'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.63%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '''3.79'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.90%  '
$ws.Range('D16').Value = '28.514.79'
$ws.Range('E16').Value = '  +5.24%  '
$ws.Range('D17').Value = '''63.38'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.78%  '
$ws.Range('D18').Value = '''232.04'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +7.55%  '
$ws.Range('D19').Value = '''7.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.55%  '
$ws.Range('E20').Value = '  +1.68%  '
$ws.Range('E21').Value = '  -0.15%  '
$ws.Range('E22').Value = '  +0.80%  '
$ws.Range('D23').Value = '''9.43'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.69%  '
$ws.Range('D24').Value = '''1.97'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.75%  '
$ws.Range('D25').Value = '''152.51'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').Value = '''15.32'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.21%  '
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  +1.32%  '
$ws.Range('E31').Value = '  +1.55%  '
$ws.Range('E32').Value = '  +1.22%  '
$ws.Range('E33').Value = '  +0.76%  '
$ws.Range('D34').Value = '1.423.92'
$ws.Range('E34').Value = '  -0.66%  '
$ws.Range('D35').Value = '''1.60'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('E36').Value = '  -3.90%  '
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').Value = '''0.0168'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.57%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').Value = '''0.544'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.62%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '''2.52'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.99%  '
$ws.Range('E41').Value = '  +2.49%  '
$ws.Range('E42').Value = '  -2.47%  '
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('D44').Value = '''0.984'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E45').Value = '  +6.60%  '
$ws.Range('D46').Value = '''64.96'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.47%  '
$ws.Range('D47').Value = '1.741.37'
$ws.Range('E47').Value = '  +2.88%  '
$ws.Range('D48').Value = '''2.15'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.73%  '
$ws.Range('D49').Value = '''87.54'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.52%  '
$ws.Range('D50').Value = '0.0₆0107'
$ws.Range('E50').Value = '  +8.35%  '
$ws.Range('E51').Value = '  +0.75%  '
